# Applies updated forecast-error values to rows 2-11 (columns B:G)
# per "Fixed naive component forecaster bug - Presentation state 11.02."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.2118531303557739, 1.305987383552482, 5.35557751699296, 2.314212072605482, 2.327425548436394, 51)
    3  = @(0.02283322487048015, 1.418515438886082, 5.16977209240434, 2.273713282805099, 2.296681442697479, 50)
    4  = @(0.2134349945989513, 1.160232958483259, 3.154048196090197, 1.775964018805054, 1.781362956999537, 49)
    5  = @(0.03943674361467967, 1.340425094488151, 5.263324471644343, 2.294193643013672, 2.318128933296117, 48)
    6  = @(0.1969686117119996, 1.570705429436083, 5.967394291932829, 2.442825063718814, 2.461194875607015, 47)
    7  = @(0.06983087975920006, 1.416739844804068, 4.44628832700042, 2.108622376576807, 2.130753392377103, 46)
    8  = @(0.2152879885089396, 1.520943876705648, 5.737038814213995, 2.395211642885446, 2.412472540257287, 45)
    9  = @(0.05959129116111418, 1.434755410722559, 5.152491520827321, 2.269910024830791, 2.295361234929747, 44)
    10 = @(0.1780234371580798, 1.510931085580928, 5.228947591369755, 2.286689220547855, 2.306729160466042, 43)
    11 = @(0.1196591527663331, 1.488698517646936, 4.992931766918498, 2.234486913570652, 2.258327489652168, 42)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
